{"js": "// Replace the double-curly-brace template placeholders with single-curly-brace\n// placeholders throughout the body (e.g. \"{{no_sk}}\" -> \"{no_sk}\"), and add a\n// trailing empty paragraph at the end of the document body.\n\nconst placeholders = [\n  \"no_sk\",\n  \"judul\",\n  \"semester\",\n  \"tanggal\",\n  \"nama_dekan\",\n  \"nip_dekan\",\n  \"ttd_base64\",\n];\n\nfor (const name of placeholders) {\n  const searchResults = context.document.body.search(\"{{\" + name + \"}}\", {\n    matchCase: true,\n  });\n  searchResults.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < searchResults.items.length; i++) {\n    searchResults.items[i].insertText(\"{\" + name + \"}\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Add a new empty paragraph at the very end of the document body.\ncontext.document.body.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# sk_wali_mhs_aktif.docx template touch-up:\n#   - collapse the Jinja-style \"{{placeholder}}\" tags down to single-brace\n#     \"{placeholder}\" tags (no_sk, judul, semester, tanggal, nama_dekan,\n#     nip_dekan, ttd_base64) so the backend's simple templating engine can\n#     substitute them.\n#   - leave a trailing blank paragraph at the end of the body (room for the\n#     signature image preview).\n\n$d = $word.ActiveDocument\n\n$placeholders = @(\n    \"no_sk\",\n    \"judul\",\n    \"semester\",\n    \"tanggal\",\n    \"nama_dekan\",\n    \"nip_dekan\",\n    \"ttd_base64\"\n)\n\nforeach ($name in $placeholders) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        \"{{\" + $name + \"}}\",   # FindText\n        $false,                # MatchCase\n        $false,                # MatchWholeWord\n        $false,                # MatchWildcards\n        $false,                # MatchSoundsLike\n        $false,                # MatchAllWordForms\n        $true,                 # Forward\n        1,                     # Wrap -> wdFindContinue\n        $false,                # Format\n        \"{\" + $name + \"}\",     # ReplaceWith\n        2                      # Replace -> wdReplaceAll\n    )\n}\n\n# Add a new empty paragraph at the very end of the document.\n$d.Content.InsertParagraphAfter()\n"}
